$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.103.55"
$ws.Range("E2").Value = "  -1.44%  "

$ws.Range("D3").Value = "2.105.13"
$ws.Range("E3").Value = "  -0.22%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.79%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "350.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.23%  "

$ws.Range("E6").Value = "  -0.77%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5167"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.64%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4483"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.58%  "

$ws.Range("E9").Value = "  -4.49%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08953"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.61%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.178"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.61%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.79"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.81%  "

$ws.Range("D13").Value = "2.110.70"
$ws.Range("E13").Value = "  -0.14%  "

$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.263"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.82%  "

$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.770"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.42%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "99.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.17%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001150"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.33%  "

$ws.Range("E18").Value = "  -0.77%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "20.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.82%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06661"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.45%  "

$ws.Range("E21").Value = "  -0.76%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.271"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.11%  "

$ws.Range("D23").Value = "30.204.53"
$ws.Range("E23").Value = "  -1.35%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.18%  "

$ws.Range("E25").Value = "  -0.41%  "

$ws.Range("D26").Value = "2.355.32"
$ws.Range("E26").Value = "  -0.33%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.19%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.553"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.34%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "163.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.30%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.77"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.08%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.185"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.04%  "

$ws.Range("E32").Value = "  -0.08%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.652"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.59%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.271"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.39%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.967"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.05%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.912"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.31%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.20"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.25%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02593"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.56%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06851"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.46%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2325"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.44%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.56"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.43%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6860"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.00%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.258"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.24%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6451"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.10%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.303"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.30%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000368"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.30%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.667"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.23%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "84.18"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.21%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.224"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.17%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07237"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.71%  "
